# filter() method for featureGroupsSet and negate support for features method
#
# This updates the "fGroups" implementation status matrix:
#  - several methods that were partially done ("X?" in column B, or an "X"
#    mark in a different status column) are now marked fully "done" (col G)
#  - the featureTable() and groupInfo() status marks move from their old
#    "almost as-is" / "done" columns over to new locations
#  - two rows (getXCMSnExp / reportCSV) get a remark in a new column H about
#    possibly waiting for the autoID branch (negate support follow-up)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 7 (as.data.table): mark done
$ws.Range("G7").Value = "X"

# row 8 (averageGroups): was "X?" (almost as-is) -> now plain "X", plus done
$ws.Range("B8").Value = "X"
$ws.Range("G8").Value = "X"

# row 12 (export): mark done
$ws.Range("G12").Value = "X"

# row 13 (featureTable): move mark from column B to column C, plus done
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "X"
$ws.Range("G13").Value = "X"

# row 14 (filter): mark done
$ws.Range("G14").Value = "X"

# row 23 (getXCMSnExp): note about waiting for autoID branch
$ws.Range("H23").Value = "maybe wait for autoID branch"

# row 24 (groupInfo): move mark from column C to column G (done)
$ws.Range("C24").ClearContents()
$ws.Range("G24").Value = "X"

# row 30 (initialize): was "X?" -> "X", plus done
$ws.Range("B30").Value = "X"
$ws.Range("G30").Value = "X"

# row 31 (length): was "X?" -> "X", plus done
$ws.Range("B31").Value = "X"
$ws.Range("G31").Value = "X"

# row 32 (names): mark done
$ws.Range("G32").Value = "X"

# row 34 (plotChord): was "X?" -> "X", plus done
$ws.Range("B34").Value = "X"
$ws.Range("G34").Value = "X"

# row 35 (plotEIC): mark done
$ws.Range("G35").Value = "X"

# row 36 (plotEICHash): mark done
$ws.Range("G36").Value = "X"

# row 37 (plotInt): mark done
$ws.Range("G37").Value = "X"

# row 38 (plotIntHash): mark done
$ws.Range("G38").Value = "X"

# row 44 (removeGroups): was "X?" -> "X", plus done
$ws.Range("B44").Value = "X"
$ws.Range("G44").Value = "X"

# row 48 (reportCSV): note about waiting for autoID branch
$ws.Range("H48").Value = "maybe wait for autoID branch"

# row 50 (reportPDF): was "X?" -> "X", plus done
$ws.Range("B50").Value = "X"
$ws.Range("G50").Value = "X"

# row 51 (updateFeatIndex): was "X?" -> "X", plus done
$ws.Range("B51").Value = "X"
$ws.Range("G51").Value = "X"

# row 53 (done/groupAlgorithm): mark done
$ws.Range("G53").Value = "X"

# keep the same selection Excel ended up with when these edits were made
$ws.Range("G15").Select()
